# Insert a new data row at row 135 (shifting existing rows 135..214 down to 136..215)
# and populate it with the new weekly price record, per the commit's intent of adding
# a new "Fruta / hortaliza, semanal" observation to the Achicoria series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 135 - this shifts all data from row 135
# downward by one row (old row 135 becomes row 136, ..., old row 214 becomes row 215).
$ws.Rows.Item(135).Insert()

# Populate the newly inserted row 135 with the new record.
$ws.Range("A135").Value = 3
$ws.Range("B135").Value = "Femacal de La Calera"
$ws.Range("C135").Value = "Coquimbo"
$ws.Range("D135").Value = 44767
$ws.Range("E135").Value = 5
$ws.Range("F135").Value = 100112010
$ws.Range("G135").Value = "Achicoria"
$ws.Range("H135").Value = "Sin especificar"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 50
$ws.Range("K135").Value = 7000
$ws.Range("L135").Value = 7000
$ws.Range("M135").Value = 7000
$ws.Range("N135").Value = "$/caja 16 unidades"
$ws.Range("O135").Value = "Provincia de Quillota"
$ws.Range("P135").Value = 438
$ws.Range("Q135").Value = 16
$ws.Range("R135").Value = "Hortaliza"
